# Apply updated crypto price/volume data to Sheet1 (inline text values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.862.47'
$ws.Range('E2').Value = '  +2.25%  '
$ws.Range('D3').Value = '2.492.34'
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.01'
$ws.Range('E5').Value = '  +3.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.55'
$ws.Range('E6').Value = '  +4.19%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  +2.97%  '
$ws.Range('D9').Value = '2.517.23'
$ws.Range('E9').Value = '  +3.00%  '
$ws.Range('E10').Value = '  +4.80%  '
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.24'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').Value = '2.946.90'
$ws.Range('E14').Value = '  +2.79%  '
$ws.Range('D15').Value = '58.838.36'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('E16').Value = '  +3.46%  '
$ws.Range('E17').Value = '  +3.47%  '
$ws.Range('D18').Value = '2.509.29'
$ws.Range('E18').Value = '  +2.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.71'
$ws.Range('E19').Value = '  +2.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '322.32'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('E22').Value = '  +8.50%  '
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.01'
$ws.Range('E24').Value = '  +4.08%  '
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.53'
$ws.Range('E28').Value = '  +4.37%  '
$ws.Range('E29').Value = '  +6.00%  '
$ws.Range('E30').Value = '  +4.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '171.85'
$ws.Range('E31').Value = '  +0.83%  '
$ws.Range('E32').Value = '  +5.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.36'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.23'
$ws.Range('E36').Value = '  +3.09%  '
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.00'
$ws.Range('E38').Value = '  +1.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.52'
$ws.Range('E39').Value = '  +4.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.75'
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.790'
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '280.83'
$ws.Range('E42').Value = '  +3.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.49'
$ws.Range('E43').Value = '  +3.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.11'
$ws.Range('E44').Value = '  +3.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '132.49'
$ws.Range('E45').Value = '  +9.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.598'
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('E47').Value = '  +2.59%  '
$ws.Range('E48').Value = '  +5.47%  '
$ws.Range('E49').Value = '  +4.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.18'
$ws.Range('E50').Value = '  +3.59%  '
$ws.Range('D51').Value = '1.759.96'
$ws.Range('E51').Value = '  +3.07%  '
